$p = $ppt.ActivePresentation

# Slide 4 ("How To Use Local Storage - Save data") - Picture 4
$s4 = $p.Slides.Item(4)
$pic4 = $s4.Shapes.Item(2)
$pic4.Left = 380.6819152832031
$pic4.Top = 185.03347778320312
$pic4.Width = 494.43212890625
$pic4.Height = 123.3488998413086

# Slide 5 ("How to Use Local Storage - Load Data") - Picture 4
$s5 = $p.Slides.Item(5)
$pic5 = $s5.Shapes.Item(2)
$pic5.Left = 404.1615905761719
$pic5.Top = 178.2234649658203
$pic5.Width = 457.29638671875
$pic5.Height = 170.89276123046875
